$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark (currently sitting
#     right after "Dr. Prakash P. Patel" at the end of the document).
#     It gets re-created further up in the body below, and bookmark names
#     must be unique, so drop the old one first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: turn "... with PASS/FAIL result." into
#     "... with XOXX result." (XOXX rendered in the Algerian font), and
#     re-plant the "_GoBack" bookmark as a collapsed range right after
#     "XOXX " and before "result.".
$rng = $d.Content
$found = $rng.Find.Execute("PASS/FAIL", $true, $false, $false, $false, $false, $true, 1, $false, "XOXX", 2)

if ($found) {
    # $rng now spans exactly the freshly-inserted "XOXX" -> give it the
    # Algerian font, matching the surrounding Algerian-styled fill-ins.
    $rng.Font.Name = "Algerian"

    $afterXoxx = $rng.End

    # Find the end of the enclosing paragraph (before its end-of-paragraph
    # mark) so we can isolate the trailing " result." text.
    $paraRng = $rng.Duplicate
    $paraRng.Expand(4) | Out-Null
    $paraTextEnd = $paraRng.End - 1

    # The paragraph still reads "...XOXX result." at this point (one
    # trailing run covering " result."). The space right after "XOXX" can
    # stay put; cut out only "result." itself (leaving the space in place)
    # so it can be retyped as its own run after the bookmark.
    $tail = $d.Range($afterXoxx + 1, $paraTextEnd)
    $tail.Delete()

    # Re-insert "result." first...
    $resultPos = $d.Range($afterXoxx + 1, $afterXoxx + 1)
    $resultPos.InsertAfter("result.")

    # ...then drop the (now collapsed) "_GoBack" bookmark right in front of
    # it, between the space and "result.".
    $bmPos = $d.Range($afterXoxx + 1, $afterXoxx + 1)
    $d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null
}
